$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "43.627.65"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -6.27%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.598.30"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.22%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "301.40"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -2.10%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "96.51"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -3.77%  "

# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  -3.84%  "

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  +0.16%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.559"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.35%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.82"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -6.61%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0816"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -3.60%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.82"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -4.17%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "2.997.09"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.26%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.15%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.599.70"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.23%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.891"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -3.40%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "14.35"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.09%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "43.648.75"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -6.37%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.69%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0979"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -3.89%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.33"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -5.42%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "72.91"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +1.86%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "265.42"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.42%  "

# Row 24
$ws.Range("B24").Value = "ImmutableX"
$ws.Range("C24").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.22"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.23%  "

# Row 25
$ws.Range("B25").Value = "PancakeSwap"
$ws.Range("C25").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.93"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.49%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "29.25"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +3.22%  "

# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.32%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.29"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -3.37%  "

# Row 29
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.22"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -3.59%  "

# Row 30
$ws.Range("B30").Value = "InjectiveProtocol"
$ws.Range("C30").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "37.58"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.97%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.07"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -5.85%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.62"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -0.76%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.21"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -1.04%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "152.14"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.10%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.26%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.0811"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -3.74%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.117"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.99%  "

# Row 38
$ws.Range("B38").Value = "EnergySwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "24.53"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +6.14%  "

# Row 39
$ws.Range("B39").Value = "Stellar"
$ws.Range("C39").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.121"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -1.76%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "16.68"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +3.57%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.88%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0314"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -5.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.86"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -5.42%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.042.75"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -3.65%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.07%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "88.04"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -5.35%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.10"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.41%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +5.12%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.854.24"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.25%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "105.95"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -3.12%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.191"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -4.63%  "

